$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old contents but keep existing cell formatting (so styles.xml stays untouched)
$ws.Range("A1:T11").ClearContents()

# --- Header row (unchanged content, rewritten to anchor shared-string indices 0-19) ---
$headers = @(
    "Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity"
)
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Sending/Ligand/Receptor/Target cluster labels (columns A-D), written column-major ---
# so the new unique shared string "Resolving-Mac" is interned right after "Inflammatory-Mac"
# (matching the order produced by the original authoring tool).
$colA = @(
    "Inflammatory-Mac", "Inflammatory-Mac", "Inflammatory-Mac", "Inflammatory-Mac", "Inflammatory-Mac", "Resolving-Mac", "Resolving-Mac", "Resolving-Mac", "Resolving-Mac", "Resolving-Mac"
)
$colB = @(
    "Ifng", "Ifng", "Ifng", "Ifng", "Ifng", "Ifng", "Ifng", "Ifng", "Ifng", "Ifng"
)
$colC = @(
    "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1", "Ifngr1"
)
$colD = @(
    "ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac", "ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac"
)
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11)
for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}

# --- Numeric columns E-T ---
$numData = @{
    2 = @(2, 0.6666666666666666, 0.4227963333333333, 1.268389, 0.9187469758736285, 0.9187469758736285, 3, 1, 46.38927966666667, 139.167839, 0.1490077759465061, 0.1506561784680905, 19.61321734904122, 176.518956141371, 0.1369004435325077, 0.1384149083642358)
    3 = @(2, 0.6666666666666666, 0.4227963333333333, 1.268389, 0.9187469758736285, 0.9187469758736285, 3, 1, 42.65870533333333, 127.976116, 0.1370247361492196, 0.1385405759713568, 18.03594419968044, 162.323497797124, 0.1258910619569773, 0.1272837352094748)
    4 = @(2, 0.6666666666666666, 0.4227963333333333, 1.268389, 0.9187469758736285, 0.9187469758736285, 3, 1, 124.6999866666666, 374.09996, 0.4005508989851949, 0.4049820040581756, 52.72269712938221, 474.5042741644399, 0.3680049271261111, 0.3720759915116904)
    5 = @(2, 0.6666666666666666, 0.4227963333333333, 1.268389, 0.9187469758736285, 0.9187469758736285, 2, 1, 10.21895, 20.4379, 0.03282445907653739, 0.02212505369083864, 4.320534590516666, 25.9232075431, 0.0301573725112564, 0.02032732616949966)
    6 = @(2, 0.6666666666666666, 0.4227963333333333, 1.268389, 0.9187469758736285, 0.9187469758736285, 3, 1, 87.35427866666667, 262.062836, 0.2805921298425419, 0.2836961878115384, 36.93306872124489, 332.397618491204, 0.2577931707467759, 0.2606450146187279)
    7 = @(1, 0.3333333333333333, 0.03739166666666666, 0.112175, 0.08125302412637153, 0.08125302412637155, 3, 1, 46.38927966666667, 139.167839, 0.1490077759465061, 0.1506561784680905, 1.734572482202778, 15.611152339825, 0.01210733241399842, 0.0122412701038547)
    8 = @(1, 0.3333333333333333, 0.03739166666666666, 0.112175, 0.08125302412637153, 0.08125302412637155, 3, 1, 42.65870533333333, 127.976116, 0.1370247361492196, 0.1385405759713568, 1.595080090255555, 14.3557208123, 0.01113367419224223, 0.01125684076188207)
    9 = @(1, 0.3333333333333333, 0.03739166666666666, 0.112175, 0.08125302412637153, 0.08125302412637155, 3, 1, 124.6999866666666, 374.09996, 0.4005508989851949, 0.4049820040581756, 4.662740334777776, 41.96466301299999, 0.03254597185908385, 0.03290601254648524)
    10 = @(1, 0.3333333333333333, 0.03739166666666666, 0.112175, 0.08125302412637153, 0.08125302412637155, 2, 1, 10.21895, 20.4379, 0.03282445907653739, 0.02212505369083864, 0.3821035720833333, 2.2926214325, 0.002667086565280987, 0.001797727521338978)
    11 = @(1, 0.3333333333333333, 0.03739166666666666, 0.112175, 0.08125302412637153, 0.08125302412637155, 3, 1, 87.35427866666667, 262.062836, 0.2805921298425419, 0.2836961878115384, 3.266322069811111, 29.3968986283, 0.02279895909576603, 0.02305117319281057)
}
foreach ($r in $rowNums) {
    $vals = $numData[$r]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($r, $j + 5).Value = $vals[$j]
    }
}
